$wb = $excel.ActiveWorkbook
$firstSheet = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($firstSheet)
$ws.Name = "arma validation"
$ws.Range("A1").Value = "mean"
$ws.Range("A2").Value = "std"
$ws.Range("A3").Value = "min"
$ws.Range("A4").Value = "max"
$ws.Range("A5").Value = "kurtosis"
$ws.Range("A6").Value = "skewness"
$ws.Range("A7").Value = "Historical"
$ws.Range("A8").Value = "Synthetic"
Write-Host "done"
